# updated anomaly detection, added XmR
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# set_2 (sheet2): trim the F:H "extra" columns, replace the B:E sample
# data, and extend the table down to 20 samples (rows 2-21).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("set_2")

# Row 1 header: drop the old F/G/H (5,6,7) headers entirely.
$ws2.Range("F1:H1").ClearContents() | Out-Null

# New B:E values for sample rows 2-21 (sample index 1-20).
$set2Data = @(
    @(44.01, 26,    24,    34),
    @(50,    48,    51,    43),
    @(32,    28,    26,    22),
    @(52,    55,    56,    44),
    @(16,    16,    21,    26),
    @(36,    36,    35,    31),
    @(21,    22,    18,    21),
    @(29,    21,    23,    22),
    @(26,    46,    44,    14),
    @(24,    22,    22,    44),
    @(18,    24,    24,    49),
    @(24,    20,    26,    23),
    @(19,    21,    27,    28),
    @(8,     11,    12,    12),
    @(24,    18.1,  27,    24),
    @(56,    52,    56,    50),
    @(32.01, 22,    18,    25),
    @(8,     12.01, 11,    17),
    @(51,    54,    52.01, 49),
    @(30,    28,    35,    22.01)
)

for ($i = 0; $i -lt $set2Data.Length; $i++) {
    $r = $i + 2
    $row = $set2Data[$i]
    $ws2.Cells.Item($r, 1).Value = $i + 1
    $ws2.Cells.Item($r, 2).Value = $row[0]
    $ws2.Cells.Item($r, 3).Value = $row[1]
    $ws2.Cells.Item($r, 4).Value = $row[2]
    $ws2.Cells.Item($r, 5).Value = $row[3]
}

# Rows 13-21 are brand new cells (rows 2-12's B:E already carried the
# "0.000" number-format style from the original file); give them the
# same "0.000" number format so they pick up the same shared style.
$ws2.Range("B13:E21").NumberFormat = "0.000"

# Old F/G/H data (rows 2-11) becomes blank, but keeps its number-format
# style; row 12 keeps only a blank F cell with that style; rows 13-21
# never had F/G/H at all.
for ($r = 2; $r -le 11; $r++) {
    $ws2.Cells.Item($r, 6).Value = ""
    $ws2.Cells.Item($r, 7).Value = ""
    $ws2.Cells.Item($r, 8).Value = ""
}
$ws2.Cells.Item(12, 6).Value = ""

# Move the active cell for set_2 to C17.
[void]$ws2.Range("C17").Select()

# ---------------------------------------------------------------------
# set_3 (sheet3): replace B2:B11 values and extend down to row 21.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("set_3")

$set3Data = @(12, 14, 16, 18, 16, 14, 12, 12, 32, 16, 18, 16, 14, 12, 16, 18, 12, 19, 18, 21)

for ($i = 0; $i -lt $set3Data.Length; $i++) {
    $r = $i + 2
    $ws3.Cells.Item($r, 1).Value = $i + 1
    $ws3.Cells.Item($r, 2).Value = $set3Data[$i]
}

# ---------------------------------------------------------------------
# set_1 (sheet1): just move the selection to C23.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("set_1")
[void]$ws1.Range("C23").Select()

# ---------------------------------------------------------------------
# Finally, move the active cell for set_3 to B10. Selecting a range
# also activates its sheet, so doing this last makes set_3 the active
# sheet/tab (tabSelected + workbook activeTab), matching the tab
# switch away from set_6.
# ---------------------------------------------------------------------
[void]$ws3.Range("B10").Select()
